$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 421.5
$ws.Range("I2").Value = 428.66666
$ws.Range("K2").Value = 428.66666
$ws.Range("M2").Value = -315.66666
$ws.Range("H9").Value = 464.7143
$ws.Range("I9").Value = 555.1818
$ws.Range("J9").Value = 133
$ws.Range("K9").Value = 555.1818
$ws.Range("L9").Value = 133
$ws.Range("M9").Value = -386.1818
$ws.Range("N9").Value = -471
$ws.Range("H44").Value = 36666.668
$ws.Range("J44").Value = 36666.668
$ws.Range("L44").Value = 36666.668
$ws.Range("N44").Value = -37590.668
$ws.Range("H86").Value = 4495.5386
$ws.Range("I86").Value = 4370.1665
$ws.Range("K86").Value = 4370.1665
$ws.Range("M86").Value = -3247.1665
$ws.Range("H89").Value = 4495.5386
$ws.Range("I89").Value = 4370.1665
$ws.Range("K89").Value = 21850.8325
$ws.Range("M89").Value = -16234.8325
$ws.Range("H111").Value = 3174.125
$ws.Range("I111").Value = 3349
$ws.Range("K111").Value = 10047
$ws.Range("M111").Value = -6980
$ws.Range("H112").Value = 2914.5557
$ws.Range("J112").Value = 2914.5557
$ws.Range("L112").Value = 8743.667099999999
$ws.Range("N112").Value = -10959.6671
$ws.Range("H113").Value = 11998.333
$ws.Range("I113").Value = 11998
$ws.Range("J113").Value = 11998.5
$ws.Range("K113").Value = 11998
$ws.Range("L113").Value = 11998.5
$ws.Range("M113").Value = -8744
$ws.Range("N113").Value = -18506.5
$ws.Range("H127").Value = 2619.8
$ws.Range("I127").Value = 2649.75
$ws.Range("J127").Value = 2500
$ws.Range("K127").Value = 7949.25
$ws.Range("L127").Value = 7500
$ws.Range("M127").Value = -2989.25
$ws.Range("N127").Value = -17420
$ws.Range("H129").Value = 2196.6
$ws.Range("I129").Value = 998
$ws.Range("J129").Value = 2995.6667
$ws.Range("K129").Value = 2994
$ws.Range("L129").Value = 8987.000100000001
$ws.Range("M129").Value = 2006
$ws.Range("N129").Value = -18987.0001
$ws.Range("H131").Value = 4848.4165
$ws.Range("J131").Value = 9598.6
$ws.Range("L131").Value = 28795.8
$ws.Range("N131").Value = -38875.8
$ws.Range("H132").Value = 41671180
$ws.Range("I132").Value = 41671180
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 125013540
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -125011010
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 13890555
$ws.Range("I137").Value = 55556870
$ws.Range("J137").Value = 1782.2778
$ws.Range("K137").Value = 166670610
$ws.Range("L137").Value = 5346.8334
$ws.Range("M137").Value = -166668060
$ws.Range("N137").Value = -10446.8334
$ws.Range("H141").Value = 3497.3125
$ws.Range("I141").Value = 3368.7273
$ws.Range("K141").Value = 10106.1819
$ws.Range("M141").Value = -4926.1819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1528.2858
$ws.Range("I45").Value = 1283
$ws.Range("K45").Value = 1283
$ws.Range("M45").Value = -906
$ws.Range("H74").Value = 2260.1365
$ws.Range("I74").Value = 1653
$ws.Range("J74").Value = 4992.25
$ws.Range("K74").Value = 1653
$ws.Range("L74").Value = 4992.25
$ws.Range("M74").Value = -779
$ws.Range("N74").Value = -6740.25
$ws.Range("H77").Value = 2260.1365
$ws.Range("I77").Value = 1653
$ws.Range("J77").Value = 4992.25
$ws.Range("K77").Value = 8265
$ws.Range("L77").Value = 24961.25
$ws.Range("M77").Value = -3897
$ws.Range("N77").Value = -33697.25
$ws.Range("H102").Value = 2632.4736
$ws.Range("I102").Value = 2444.2144
$ws.Range("J102").Value = 3159.6
$ws.Range("K102").Value = 2444.2144
$ws.Range("L102").Value = 3159.6
$ws.Range("M102").Value = -822.2143999999998
$ws.Range("N102").Value = -6403.6
$ws.Range("H110").Value = 1632.6666
$ws.Range("I110").Value = 1632.6666
$ws.Range("K110").Value = 1632.6666
$ws.Range("M110").Value = 412.3334
$ws.Range("H132").Value = 4175.7095
$ws.Range("I132").Value = 3312.8518
$ws.Range("K132").Value = 9938.555399999999
$ws.Range("M132").Value = -7408.555399999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 31200
$ws.Range("I45").Value = 26500
$ws.Range("K45").Value = 26500
$ws.Range("M45").Value = -25907
$ws.Range("H50").Value = 19713.818
$ws.Range("J50").Value = 32000
$ws.Range("L50").Value = 32000
$ws.Range("N50").Value = -33250
$ws.Range("H58").Value = 3025.5454
$ws.Range("I58").Value = 3025.5454
$ws.Range("K58").Value = 3025.5454
$ws.Range("M58").Value = -2822.5454
$ws.Range("H86").Value = 12669.077
$ws.Range("I86").Value = 11101
$ws.Range("J86").Value = 16197.25
$ws.Range("K86").Value = 11101
$ws.Range("L86").Value = 16197.25
$ws.Range("M86").Value = -9978
$ws.Range("N86").Value = -18443.25
$ws.Range("H89").Value = 12669.077
$ws.Range("I89").Value = 11101
$ws.Range("J89").Value = 16197.25
$ws.Range("K89").Value = 55505
$ws.Range("L89").Value = 80986.25
$ws.Range("M89").Value = -49889
$ws.Range("N89").Value = -92218.25
$ws.Range("H96").Value = 22086.5
$ws.Range("J96").Value = 22086.5
$ws.Range("L96").Value = 22086.5
$ws.Range("N96").Value = -27578.5
$ws.Range("H105").Value = 499
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("M122").Value = -3548.5
$ws.Range("H136").Value = 3025.5454
$ws.Range("I136").Value = 3025.5454
$ws.Range("K136").Value = 9076.636200000001
$ws.Range("M136").Value = -6526.636200000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3640174.2
$ws.Range("I4").Value = 1992685.2
$ws.Range("K4").Value = 5978055.6
$ws.Range("M4").Value = -5977943.6
$ws.Range("H16").Value = 632.6667
$ws.Range("J16").Value = 699
$ws.Range("L16").Value = 2097
$ws.Range("N16").Value = -2443
$ws.Range("H34").Value = 49881.273
$ws.Range("J34").Value = 54816.95
$ws.Range("L34").Value = 164450.85
$ws.Range("N34").Value = -164618.85
$ws.Range("H69").Value = 2136.1924
$ws.Range("I69").Value = 1878.6
$ws.Range("J69").Value = 2197.524
$ws.Range("K69").Value = 5635.799999999999
$ws.Range("L69").Value = 6592.572
$ws.Range("M69").Value = -4824.799999999999
$ws.Range("N69").Value = -8214.572
$ws.Range("H72").Value = 2136.1924
$ws.Range("I72").Value = 1878.6
$ws.Range("J72").Value = 2197.524
$ws.Range("K72").Value = 16907.4
$ws.Range("L72").Value = 19777.716
$ws.Range("M72").Value = -12851.4
$ws.Range("N72").Value = -27889.716
$ws.Range("H87").Value = 804.6667
$ws.Range("I87").Value = 804.6667
$ws.Range("K87").Value = 2414.0001
$ws.Range("M87").Value = -1166.0001
$ws.Range("H90").Value = 804.6667
$ws.Range("I90").Value = 804.6667
$ws.Range("K90").Value = 7242.0003
$ws.Range("M90").Value = -1002.0003
$ws.Range("H107").Value = 427.35483
$ws.Range("J107").Value = 416
$ws.Range("L107").Value = 1248
$ws.Range("N107").Value = -5088
$ws.Range("H113").Value = 2235.7
$ws.Range("J113").Value = 2298
$ws.Range("L113").Value = 6894
$ws.Range("N113").Value = -11234
$ws.Range("H129").Value = 1477
$ws.Range("I129").Value = 796.8333
$ws.Range("J129").Value = 2060
$ws.Range("K129").Value = 2390.4999
$ws.Range("L129").Value = 6180
$ws.Range("M129").Value = 2609.5001
$ws.Range("N129").Value = -16180
$ws.Range("H131").Value = 1376.5238
$ws.Range("I131").Value = 1126.6666
$ws.Range("J131").Value = 1476.4667
$ws.Range("K131").Value = 3379.9998
$ws.Range("L131").Value = 4429.4001
$ws.Range("M131").Value = 1660.0002
$ws.Range("N131").Value = -14509.4001
$ws.Range("H132").Value = 900
$ws.Range("I132").Value = 914
$ws.Range("K132").Value = 8226
$ws.Range("M132").Value = -5696
$ws.Range("H139").Value = 2794.875
$ws.Range("I139").Value = 1075.1428
$ws.Range("K139").Value = 3225.4284
$ws.Range("M139").Value = 1914.5716
$ws.Range("H140").Value = 385861.06
$ws.Range("I140").Value = 385861.06
$ws.Range("K140").Value = 1157583.18
$ws.Range("M140").Value = -1152403.18

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11497637
$ws.Range("J132").Value = 33336578
$ws.Range("L132").Value = 100009734
$ws.Range("N132").Value = -100014794

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 55046.145
$ws.Range("J98").Value = 55046.145
$ws.Range("L98").Value = 55046.145
$ws.Range("N98").Value = -61036.145
$ws.Range("H132").Value = 2820.9524
$ws.Range("I132").Value = 3768.6
$ws.Range("J132").Value = 2524.8125
$ws.Range("K132").Value = 11305.8
$ws.Range("L132").Value = 7574.4375
$ws.Range("M132").Value = -8775.799999999999
$ws.Range("N132").Value = -12634.4375
$ws.Range("H136").Value = 250001500
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H100").Value = 1920.2
$ws.Range("I100").Value = 1234
$ws.Range("K100").Value = 2468
$ws.Range("M100").Value = -1927
$ws.Range("H107").Value = 2600.5
$ws.Range("J107").Value = 2202
$ws.Range("L107").Value = 6606
$ws.Range("N107").Value = -10446
$ws.Range("H136").Value = 9217.786
$ws.Range("I136").Value = 9854.23
$ws.Range("K136").Value = 29562.69
$ws.Range("M136").Value = -27012.69

